$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Move the "common/WhirlyGlobeLib/" link from J24 to J25 (keep formatting)
# ---------------------------------------------------------------------------
$ws.Range("J24").Copy()
$ws.Range("J25").PasteSpecial(-4122)
$ws.Range("J25").Value = $ws.Range("J24").Value2
$ws.Range("J24").Clear()

# ---------------------------------------------------------------------------
# 2. Add the new "common WhirlyGlobeLib" summary block (rows 46-66), mirroring
#    the style of the existing "common WhirlyGlobeLib" block at rows 24-44.
# ---------------------------------------------------------------------------

# Row 46: section header (like row 24) -- copy full-row formatting first
$ws.Range("A24:J24").Copy()
$ws.Range("A46:J46").PasteSpecial(-4122)
$ws.Range("A46").Value = "Feb 19 - 3.0"
$ws.Range("B46").Value = "common WhirlyGlobeLib"

# Row 47: " C Header" total line, with the folder link in column J
$ws.Range("B47").Value = " C Header"
$ws.Range("E47").Value = 5588
$ws.Range("J24").Copy()
$ws.Range("J47").PasteSpecial(-4122)
$ws.Range("J47").Value = "common/WhirlyGlobeLib/"

# Row 48: " C++ " total line
$ws.Range("B48").Value = " C++ "
$ws.Range("E48").Value = 21898

# Row 49: Total line (bold, like row 27/34/41/43)
$ws.Range("B34").Copy()
$ws.Range("B49").PasteSpecial(-4122)
$ws.Range("B49").Value = "Total"
$ws.Range("E49").Value = 27486

# Row 51: "iOS WhirlyGlobeLib" sub-section header (like row 29)
$ws.Range("A29:J29").Copy()
$ws.Range("A51:J51").PasteSpecial(-4122)
$ws.Range("B51").Value = "iOS WhirlyGlobeLib"
$ws.Range("J51").Value = "ios/library/WhirlyGlobeLib/"

# Row 52-55: per-language code totals
$ws.Range("B52").Value = "C Header "
$ws.Range("E52").Value = 593
$ws.Range("B53").Value = "C++"
$ws.Range("E53").Value = 57
$ws.Range("B54").Value = "Objective C"
$ws.Range("E54").Value = 69
$ws.Range("B55").Value = "Objective C++"
$ws.Range("E55").Value = 3009

# Row 56: Total line
$ws.Range("B34").Copy()
$ws.Range("B56").PasteSpecial(-4122)
$ws.Range("B56").Value = "Total"
$ws.Range("E56").Value = 3728

# Row 58: Android sub-section header (like row 36; B36/J36 reuse the "iOS"
# labelled strings verbatim, matching a pre-existing quirk in the sheet)
$ws.Range("A36:J36").Copy()
$ws.Range("A58:J58").PasteSpecial(-4122)
$ws.Range("B58").Value = $ws.Range("B36").Value2
$ws.Range("J58").Value = $ws.Range("J36").Value2

# Row 59-62: per-language code totals
$ws.Range("B59").Value = "C Header "
$ws.Range("E59").Value = 4122
$ws.Range("B60").Value = "C++"
$ws.Range("E60").Value = 1231
$ws.Range("B61").Value = "Objective C"
$ws.Range("E61").Value = 2480
$ws.Range("B62").Value = "Objective C++"
$ws.Range("E62").Value = 20085

# Row 63: Total line
$ws.Range("B34").Copy()
$ws.Range("B63").PasteSpecial(-4122)
$ws.Range("B63").Value = "Total"
$ws.Range("E63").Value = 27955

# Row 65: grand Total line, summing the three totals above
$ws.Range("B34").Copy()
$ws.Range("B65").PasteSpecial(-4122)
$ws.Range("B65").Value = "Total"
$ws.Range("E65").Formula = "=E49+E56+E63"

# Row 66: Core Ratio line (percentage format) + "Get this higher!" note
$ws.Range("B44").Copy()
$ws.Range("B66").PasteSpecial(-4122)
$ws.Range("B66").Value = "Core Ratio"
$ws.Range("E44").Copy()
$ws.Range("E66").PasteSpecial(-4122)
$ws.Range("E66").Formula = "=E49/E65"

$ws.Range("J44").Copy()
$ws.Range("J66").PasteSpecial(-4122)
$ws.Range("J66").Value = $ws.Range("J44").Value2
$ws.Range("J44").Clear()

# ---------------------------------------------------------------------------
# 3. Update the view state to match (scrolled down, new selection)
# ---------------------------------------------------------------------------
$ws.Range("A7").Select()
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("H44").Select()
